$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(62, 8).Value2 = 2417.6667
$ws.Cells.Item(62, 10).Value2 = 2500
$ws.Cells.Item(62, 12).Value2 = 2500
$ws.Cells.Item(62, 14).Value2 = -3748
$ws.Cells.Item(65, 8).Value2 = 2417.6667
$ws.Cells.Item(65, 10).Value2 = 2500
$ws.Cells.Item(65, 12).Value2 = 12500
$ws.Cells.Item(65, 14).Value2 = -18740
$ws.Cells.Item(137, 8).Value2 = 1079
$ws.Cells.Item(137, 9).Value2 = 694.8
$ws.Cells.Item(137, 11).Value2 = 2084.4
$ws.Cells.Item(137, 13).Value2 = 465.6000000000004
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(19, 8).Value2 = 0
$ws.Cells.Item(19, 9).Value2 = 0
$ws.Cells.Item(19, 11).Value2 = 0
$ws.Cells.Item(19, 13).ClearContents()
$ws.Cells.Item(32, 8).Value2 = 3987.7778
$ws.Cells.Item(32, 9).Value2 = 3987.7778
$ws.Cells.Item(32, 11).Value2 = 3987.7778
$ws.Cells.Item(32, 13).Value2 = -3700.7778
$ws.Cells.Item(44, 8).Value2 = 31000
$ws.Cells.Item(44, 10).Value2 = 31000
$ws.Cells.Item(44, 12).Value2 = 31000
$ws.Cells.Item(44, 14).Value2 = -31976
$ws.Cells.Item(55, 8).Value2 = 20566.666
$ws.Cells.Item(55, 10).Value2 = 31750
$ws.Cells.Item(55, 12).Value2 = 31750
$ws.Cells.Item(55, 14).Value2 = -32380
$ws.Cells.Item(63, 8).Value2 = 5197.75
$ws.Cells.Item(66, 8).Value2 = 5197.75
$ws.Cells.Item(125, 8).Value2 = 22166.334
$ws.Cells.Item(125, 10).Value2 = 22166.334
$ws.Cells.Item(125, 12).Value2 = 22166.334
$ws.Cells.Item(125, 14).Value2 = -32006.334
$ws.Cells.Item(134, 8).Value2 = 107466.86
$ws.Cells.Item(134, 10).Value2 = 107466.86
$ws.Cells.Item(134, 12).Value2 = 107466.86
$ws.Cells.Item(134, 14).Value2 = -117606.86
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(35, 8).Value2 = 32026
$ws.Cells.Item(35, 9).Value2 = 0
$ws.Cells.Item(35, 11).Value2 = 0
$ws.Cells.Item(35, 13).ClearContents()
$ws.Cells.Item(82, 8).Value2 = 45000
$ws.Cells.Item(82, 9).Value2 = 0
$ws.Cells.Item(82, 11).Value2 = 0
$ws.Cells.Item(82, 13).ClearContents()
$ws.Cells.Item(85, 8).Value2 = 45000
$ws.Cells.Item(85, 9).Value2 = 0
$ws.Cells.Item(85, 11).Value2 = 0
$ws.Cells.Item(85, 13).ClearContents()
$ws.Cells.Item(128, 8).Value2 = 2700
$ws.Cells.Item(128, 9).Value2 = 2700
$ws.Cells.Item(128, 11).Value2 = 8100
$ws.Cells.Item(128, 13).Value2 = -5610
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(41, 8).Value2 = 22749.5
$ws.Cells.Item(41, 10).Value2 = 21300
$ws.Cells.Item(41, 12).Value2 = 21300
$ws.Cells.Item(41, 14).Value2 = -22156
$ws.Cells.Item(51, 8).Value2 = 20333.334
$ws.Cells.Item(51, 10).Value2 = 26500
$ws.Cells.Item(51, 12).Value2 = 26500
$ws.Cells.Item(51, 14).Value2 = -27972
$ws.Cells.Item(59, 8).Value2 = 40666.668
$ws.Cells.Item(59, 10).Value2 = 41000
$ws.Cells.Item(59, 12).Value2 = 41000
$ws.Cells.Item(59, 14).Value2 = -43290
$ws.Cells.Item(60, 8).Value2 = 25600
$ws.Cells.Item(60, 10).Value2 = 25600
$ws.Cells.Item(60, 12).Value2 = 25600
$ws.Cells.Item(60, 14).Value2 = -26622
$ws.Cells.Item(61, 8).Value2 = 20333.334
$ws.Cells.Item(61, 10).Value2 = 26500
$ws.Cells.Item(61, 12).Value2 = 26500
$ws.Cells.Item(61, 14).Value2 = -27196
$ws.Cells.Item(68, 8).Value2 = 41362.4
$ws.Cells.Item(68, 10).Value2 = 41362.4
$ws.Cells.Item(68, 12).Value2 = 41362.4
$ws.Cells.Item(68, 14).Value2 = -42860.4
$ws.Cells.Item(71, 8).Value2 = 41362.4
$ws.Cells.Item(71, 10).Value2 = 41362.4
$ws.Cells.Item(71, 12).Value2 = 124087.2
$ws.Cells.Item(71, 14).Value2 = -131575.2
$ws.Cells.Item(74, 8).Value2 = 40566.5
$ws.Cells.Item(74, 10).Value2 = 40566.5
$ws.Cells.Item(74, 12).Value2 = 40566.5
$ws.Cells.Item(74, 14).Value2 = -42314.5
$ws.Cells.Item(77, 8).Value2 = 40566.5
$ws.Cells.Item(77, 10).Value2 = 40566.5
$ws.Cells.Item(77, 12).Value2 = 121699.5
$ws.Cells.Item(77, 14).Value2 = -130435.5
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(3, 8).Value2 = 12148.818
$ws.Cells.Item(3, 9).Value2 = 9204.625
$ws.Cells.Item(3, 11).Value2 = 27613.875
$ws.Cells.Item(3, 13).Value2 = -27501.875
$ws.Cells.Item(5, 8).Value2 = 1478.4
$ws.Cells.Item(5, 9).Value2 = 1475.1538
$ws.Cells.Item(5, 10).Value2 = 1499.5
$ws.Cells.Item(5, 11).Value2 = 4425.4614
$ws.Cells.Item(5, 12).Value2 = 4498.5
$ws.Cells.Item(5, 13).Value2 = -4313.4614
$ws.Cells.Item(5, 14).Value2 = -4722.5
$ws.Cells.Item(80, 8).Value2 = 7333
$ws.Cells.Item(80, 9).Value2 = 7333
$ws.Cells.Item(80, 11).Value2 = 21999
$ws.Cells.Item(80, 13).Value2 = -21063
$ws.Cells.Item(83, 8).Value2 = 7333
$ws.Cells.Item(83, 9).Value2 = 7333
$ws.Cells.Item(83, 11).Value2 = 65997
$ws.Cells.Item(83, 13).Value2 = -61317
$ws.Cells.Item(132, 8).Value2 = 2093.6667
$ws.Cells.Item(132, 9).Value2 = 1174.3334
$ws.Cells.Item(132, 11).Value2 = 10569.0006
$ws.Cells.Item(132, 13).Value2 = -8039.000599999999
$ws.Cells.Item(133, 8).Value2 = 2416.6667
$ws.Cells.Item(133, 9).Value2 = 2416.6667
$ws.Cells.Item(133, 11).Value2 = 7250.000100000001
$ws.Cells.Item(133, 13).Value2 = -2190.000100000001
$ws.Cells.Item(135, 8).Value2 = 1478.4
$ws.Cells.Item(135, 9).Value2 = 1475.1538
$ws.Cells.Item(135, 10).Value2 = 1499.5
$ws.Cells.Item(135, 11).Value2 = 13276.3842
$ws.Cells.Item(135, 12).Value2 = 13495.5
$ws.Cells.Item(135, 13).Value2 = -10741.3842
$ws.Cells.Item(135, 14).Value2 = -18565.5
$ws.Cells.Item(137, 8).Value2 = 1186.5294
$ws.Cells.Item(137, 9).Value2 = 1135.6875
$ws.Cells.Item(137, 10).Value2 = 2000
$ws.Cells.Item(137, 11).Value2 = 3407.0625
$ws.Cells.Item(137, 12).Value2 = 6000
$ws.Cells.Item(137, 13).Value2 = 1692.9375
$ws.Cells.Item(137, 14).Value2 = -16200
$ws.Cells.Item(140, 8).Value2 = 1927.8
$ws.Cells.Item(140, 9).Value2 = 1576.4166
$ws.Cells.Item(140, 11).Value2 = 4729.2498
$ws.Cells.Item(140, 13).Value2 = 450.7502000000004
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(46, 8).Value2 = 33521.145
$ws.Cells.Item(46, 10).Value2 = 40882.668
$ws.Cells.Item(46, 12).Value2 = 40882.668
$ws.Cells.Item(46, 14).Value2 = -41194.668
$ws.Cells.Item(80, 8).Value2 = 5407.6665
$ws.Cells.Item(80, 9).Value2 = 4498.6665
$ws.Cells.Item(80, 11).Value2 = 4498.6665
$ws.Cells.Item(80, 13).Value2 = -3500.6665
$ws.Cells.Item(83, 8).Value2 = 5407.6665
$ws.Cells.Item(83, 9).Value2 = 4498.6665
$ws.Cells.Item(83, 11).Value2 = 22493.3325
$ws.Cells.Item(83, 13).Value2 = -17501.3325
$ws.Cells.Item(119, 8).Value2 = 0
$ws.Cells.Item(119, 10).Value2 = 0
$ws.Cells.Item(119, 12).Value2 = 0
$ws.Cells.Item(119, 14).ClearContents()
$ws.Cells.Item(134, 8).Value2 = 58326.332
$ws.Cells.Item(134, 10).Value2 = 58326.332
$ws.Cells.Item(134, 12).Value2 = 174978.996
$ws.Cells.Item(134, 14).Value2 = -180048.996
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value2 = 1536.5385
$ws.Cells.Item(16, 9).Value2 = 1424.5
$ws.Cells.Item(16, 10).Value2 = 1715.8
$ws.Cells.Item(16, 11).Value2 = 1424.5
$ws.Cells.Item(16, 12).Value2 = 1715.8
$ws.Cells.Item(16, 13).Value2 = -1254.5
$ws.Cells.Item(16, 14).Value2 = -2055.8
$ws.Cells.Item(93, 8).Value2 = 3499.8572
$ws.Cells.Item(93, 9).Value2 = 4399.6665
$ws.Cells.Item(93, 10).Value2 = 2825
$ws.Cells.Item(93, 11).Value2 = 4399.6665
$ws.Cells.Item(93, 12).Value2 = 2825
$ws.Cells.Item(93, 13).Value2 = -3151.6665
$ws.Cells.Item(93, 14).Value2 = -5321
$ws.Cells.Item(122, 8).Value2 = 3486.4243
$ws.Cells.Item(122, 9).Value2 = 3424.3333
$ws.Cells.Item(122, 11).Value2 = 10272.9999
$ws.Cells.Item(122, 13).Value2 = -7822.999899999999
$ws.Cells.Item(132, 8).Value2 = 8644.704
$ws.Cells.Item(132, 9).Value2 = 8540.359
$ws.Cells.Item(132, 11).Value2 = 25621.077
$ws.Cells.Item(132, 13).Value2 = -23091.077
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(54, 8).Value2 = 23918.334
$ws.Cells.Item(54, 9).Value2 = 20877.5
$ws.Cells.Item(54, 10).Value2 = 30000
$ws.Cells.Item(54, 11).Value2 = 20877.5
$ws.Cells.Item(54, 12).Value2 = 30000
$ws.Cells.Item(54, 13).Value2 = -20357.5
$ws.Cells.Item(54, 14).Value2 = -31040
$ws.Cells.Item(122, 8).Value2 = 3298.5833
$ws.Cells.Item(122, 9).Value2 = 1881.4445
$ws.Cells.Item(122, 10).Value2 = 7550
$ws.Cells.Item(122, 11).Value2 = 5644.333500000001
$ws.Cells.Item(122, 12).Value2 = 22650
$ws.Cells.Item(122, 13).Value2 = -3194.333500000001
$ws.Cells.Item(122, 14).Value2 = -27550
